$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text (string) type,
# matching the source data which stores these as text, not numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '39.596.02'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.153.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.94'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.389'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0843'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.85'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.475.01'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.72'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.802'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.46'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.144.32'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '39.593.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.48'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.05'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0842'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.62'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.65'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.76'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.56%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.56'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.67'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0616'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.87'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.09%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.63'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.90%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +22.18%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.54'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.64'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.514.36'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0916'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.77'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '49.74'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +8.35%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.07%  '
